$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Status text change: "Ready for handoff" -> "In Translation"
#    This shared string is used as the "Status" value on every sheet:
#      Overview : E2 (zh-cn status), F2 (de-de status)
#      zh-cn    : C2 (Status column)
#      de-de    : C2 (Status column)
# ---------------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"

# ---------------------------------------------------------------------------
# 2. Narrow the "Latest HO Xliff Generate Date" / date-ish columns from
#    ~17.22 chars down to ~13.41 chars:
#      Overview : columns E and F
#      zh-cn    : column C
#      de-de    : column C
# ---------------------------------------------------------------------------
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5

$zhcn.Columns.Item(3).ColumnWidth = 12.5

$dede.Columns.Item(3).ColumnWidth = 12.5
